$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the stray "x?" / "x" markers in the Professor/Lehrbeauftragter
#        columns for the first four data rows (these cells become blank,
#        keeping their existing border style). ---
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("D16").Value = ""

# --- 2. Append two new three-row blocks (Maintain/Query/Report) describing
#        "data on services provided to another faculty" (rows 18-20) and
#        "data on services used from another faculty" (rows 21-23), mirroring
#        the layout/format of the existing blocks above. ---

# Copy the formatting of the first block (rows 3-5) down onto the two new
# blocks so border styles match exactly.
$ws.Range("B3:G5").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 18 (Maintain)
$ws.Range("B18").Value = "data on services provided to another faculty"
$ws.Range("C18").Value = "Maintain"
$ws.Range("F18").Value = "x"

# Row 19 (Query)
$ws.Range("C19").Value = "Query"
$ws.Range("F19").Value = "x"

# Row 20 (Report)
$ws.Range("C20").Value = "Report"
$ws.Range("D20").Value = "x"
$ws.Range("F20").Value = "x"
$ws.Range("G20").Value = "x"

# Row 21 (Maintain)
$ws.Range("B21").Value = "data on services used from another faculty"
$ws.Range("C21").Value = "Maintain"
$ws.Range("F21").Value = "x"

# Row 22 (Query)
$ws.Range("C22").Value = "Query"
$ws.Range("F22").Value = "x"

# Row 23 (Report)
$ws.Range("C23").Value = "Report"
$ws.Range("D23").Value = "x"
$ws.Range("F23").Value = "x"
$ws.Range("G23").Value = "x"

# --- 3. Restore the selection to match the author's final cursor position. ---
$ws.Range("B21").Select()
